# Updates crypto price (D) and 1h volume change (E) columns
# to the latest refreshed values, mirroring the automated
# "Updated cryptos list" GitHub Actions data refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.392.61"
$ws.Range("D3").Value = "'1.824.80"
$ws.Range("E3").Value = '  +1.78%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = "'313.96"
$ws.Range("E5").Value = '  +1.57%  '
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").Value = "'0.4660"
$ws.Range("E7").Value = '  +4.78%  '
$ws.Range("D8").Value = "'0.3789"
$ws.Range("E8").Value = '  +3.63%  '
$ws.Range("D9").Value = "'0.07446"
$ws.Range("E9").Value = '  +2.23%  '
$ws.Range("D10").Value = "'0.8759"
$ws.Range("E10").Value = '  +2.69%  '
$ws.Range("D11").Value = "'20.81"
$ws.Range("E11").Value = '  +1.26%  '
$ws.Range("D12").Value = "'1.824.84"
$ws.Range("E12").Value = '  -3.01%  '
$ws.Range("D13").Value = "'6.687"
$ws.Range("E13").Value = '  +1.28%  '
$ws.Range("D14").Value = "'5.424"
$ws.Range("E14").Value = '  +3.20%  '
$ws.Range("D15").Value = "'92.94"
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").Value = "'0.07091"
$ws.Range("E16").Value = '  +0.33%  '
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = "'0.000008796"
$ws.Range("E18").Value = '  +1.67%  '
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("E20").Value = '  +1.76%  '
$ws.Range("D21").Value = "'27.399.99"
$ws.Range("E21").Value = '  +2.63%  '
$ws.Range("E22").Value = '  +3.66%  '
$ws.Range("D23").Value = "'10.98"
$ws.Range("E23").Value = '  +2.17%  '
$ws.Range("D24").Value = "'2.050.92"
$ws.Range("E24").Value = '  -3.79%  '
$ws.Range("E25").Value = '  -2.12%  '
$ws.Range("D26").Value = "'151.31"
$ws.Range("E26").Value = '  -0.32%  '
$ws.Range("D27").Value = "'2.256"
$ws.Range("E27").Value = '  +3.87%  '
$ws.Range("D28").Value = "'18.62"
$ws.Range("E28").Value = '  +1.53%  '
$ws.Range("D29").Value = "'5.336"
$ws.Range("E29").Value = '  +3.32%  '
$ws.Range("D30").Value = "'117.30"
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("D31").Value = "'0.08955"
$ws.Range("E31").Value = '  +2.03%  '
$ws.Range("D32").Value = "'0.7877"
$ws.Range("E32").Value = '  +6.52%  '
$ws.Range("E33").Value = '  +3.46%  '
$ws.Range("E34").Value = '  +2.30%  '
$ws.Range("D35").Value = "'2.946"
$ws.Range("E35").Value = '  +0.74%  '
$ws.Range("D36").Value = "'1.0000"
$ws.Range("E36").Value = '  -0.15%  '
$ws.Range("D37").Value = "'1.101"
$ws.Range("E37").Value = '  +1.60%  '
$ws.Range("E38").Value = '  +1.20%  '
$ws.Range("D39").Value = "'0.05255"
$ws.Range("E39").Value = '  +1.96%  '
$ws.Range("D40").Value = "'7.299"
$ws.Range("E40").Value = '  +4.26%  '
$ws.Range("D41").Value = "'0.5376"
$ws.Range("E41").Value = '  +1.95%  '
$ws.Range("D42").Value = "'2.901"
$ws.Range("E42").Value = '  +2.46%  '
$ws.Range("D43").Value = "'2.354"
$ws.Range("E43").Value = '  +20.78%  '
$ws.Range("D44").Value = "'0.1702"
$ws.Range("E44").Value = '  +1.72%  '
$ws.Range("E45").Value = '  +3.33%  '
$ws.Range("D46").Value = "'0.5106"
$ws.Range("E46").Value = '  +0.48%  '
$ws.Range("D47").Value = "'10.58"
$ws.Range("E47").Value = '  +1.72%  '
$ws.Range("D48").Value = "'105.74"
$ws.Range("E48").Value = '  +0.53%  '
$ws.Range("D49").Value = "'1.682"
$ws.Range("E49").Value = '  +1.89%  '
$ws.Range("D50").Value = "'1.000"
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("D51").Value = "'0.06384"
$ws.Range("E51").Value = '  +1.47%  '
